$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update PrefabName values in column J (sharedStrings reshuffle effect) ---
# Row 6 (Knight -> Wizard), Row 8 (BaoYu -> Bat), Row 9 (BaoYu -> BlackBoar), Row 7 (BaoYu -> Monkey)
# Order below matters: it reproduces the original author's shared-string insertion order.
$ws.Range("J6").Value = "Wizard"
$ws.Range("J8").Value = "Bat"
$ws.Range("J9").Value = "BlackBoar"
$ws.Range("J7").Value = "Monkey"

# --- Sheet view: update the active selection to J7 ---
$ws.Range("J7").Select()

$wb.Save()
